$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

function Set-PlainCell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value2 = $val
}

# Row 150
Set-PlainCell 150 1 "AMCR"
Set-TextCell 150 2 "04/06/2020"
Set-PlainCell 150 3 "Nicholas T. Long"
$ws.Cells.Item(150,4).Value2 = 240
Set-TextCell 150 5 "`$8.20"
$ws.Cells.Item(150,6).Value2 = 1968
$ws.Cells.Item(150,7).Value2 = 5.9

# Row 151
Set-PlainCell 151 1 "AAL"
Set-TextCell 151 2 "02/28/2020"
Set-PlainCell 151 3 "MICHAEL J EMBLER"
$ws.Cells.Item(151,4).Value2 = 4000
Set-TextCell 151 5 "`$19.33"
$ws.Cells.Item(151,6).Value2 = 77312
$ws.Cells.Item(151,7).Value2 = 12.92

# Row 152
Set-PlainCell 152 1 "AAL"
Set-TextCell 152 2 "02/28/2020"
Set-PlainCell 152 3 "JOHN T CAHILL"
$ws.Cells.Item(152,4).Value2 = 25000
Set-TextCell 152 5 "`$18.96"
$ws.Cells.Item(152,6).Value2 = 474125
$ws.Cells.Item(152,7).Value2 = 100

# Row 153
Set-PlainCell 153 1 "AAL"
Set-TextCell 153 2 "02/24/2020"
Set-PlainCell 153 3 "JOHN T CAHILL"
$ws.Cells.Item(153,4).Value2 = 25000
Set-TextCell 153 5 "`$25.14"
$ws.Cells.Item(153,6).Value2 = 628377.5
$ws.Cells.Item(153,7).Value2 = 100

# Row 154
Set-PlainCell 154 1 "AXP"
Set-TextCell 154 2 "10/24/2023"
Set-PlainCell 154 3 "Walter Joseph III Clayton"
$ws.Cells.Item(154,4).Value2 = 1000
Set-TextCell 154 5 "`$143.93"
$ws.Cells.Item(154,6).Value2 = 143930
$ws.Cells.Item(154,7).Value2 = 100

# Row 155
Set-PlainCell 155 1 "AXP"
Set-TextCell 155 2 "11/08/2022"
Set-PlainCell 155 3 "Walter Joseph III Clayton"
$ws.Cells.Item(155,4).Value2 = 1000
Set-TextCell 155 5 "`$149.27"
$ws.Cells.Item(155,6).Value2 = 149270
$ws.Cells.Item(155,7).Value2 = 100

# Row 156
Set-PlainCell 156 1 "AXP"
Set-TextCell 156 2 "02/11/2020"
Set-PlainCell 156 3 "Lynn Ann Pike"
$ws.Cells.Item(156,4).Value2 = 1000
Set-TextCell 156 5 "`$132.87"
$ws.Cells.Item(156,6).Value2 = 132870
$ws.Cells.Item(156,7).Value2 = 1538.46

# Row 157
Set-PlainCell 157 1 "AIG"
Set-TextCell 157 2 "06/28/2024"
Set-PlainCell 157 3 "John C Inglis"
$ws.Cells.Item(157,4).Value2 = 6.9116
Set-TextCell 157 5 "`$74.89"
$ws.Cells.Item(157,6).Value2 = 517.60350356
$ws.Cells.Item(157,7).Value2 = 1.05

# Row 158
Set-PlainCell 158 1 "AIG"
Set-TextCell 158 2 "03/14/2024"
Set-PlainCell 158 3 "John C Inglis"
$ws.Cells.Item(158,4).Value2 = 659
Set-TextCell 158 5 "`$75.39"
$ws.Cells.Item(158,6).Value2 = 49682.01
$ws.Cells.Item(158,7).Value2 = 100

# Row 159
Set-PlainCell 159 1 "AIG"
Set-TextCell 159 2 "10/30/2020"
Set-PlainCell 159 3 "INTERNATIONAL GROUP INC AMERICAN"
$ws.Cells.Item(159,4).Value2 = 625000
Set-TextCell 159 5 "`$16.00"
$ws.Cells.Item(159,6).Value2 = 10000000
$ws.Cells.Item(159,7).Value2 = 100

# Row 160
Set-PlainCell 160 1 "AIG"
Set-TextCell 160 2 "05/06/2020"
Set-PlainCell 160 3 "WILLIAM G JURGENSEN"
$ws.Cells.Item(160,4).Value2 = 20000
Set-TextCell 160 5 "`$24.10"
$ws.Cells.Item(160,6).Value2 = 481965.9999999999
$ws.Cells.Item(160,7).Value2 = 133.33

# Row 161
Set-PlainCell 161 1 "AWK"
Set-TextCell 161 2 "11/06/2023"
Set-PlainCell 161 3 "Michael Marberry"
$ws.Cells.Item(161,4).Value2 = 3786
Set-TextCell 161 5 "`$130.20"
$ws.Cells.Item(161,6).Value2 = 492937.2
$ws.Cells.Item(161,7).Value2 = 103.08

# Row 162
Set-PlainCell 162 1 "AWK"
Set-TextCell 162 2 "08/28/2023"
Set-PlainCell 162 3 "MARTHA CLARK GOSS"
$ws.Cells.Item(162,4).Value2 = 56
Set-TextCell 162 5 "`$141.07"
$ws.Cells.Item(162,6).Value2 = 7899.92
$ws.Cells.Item(162,7).Value2 = 0.16

# Row 163
Set-PlainCell 163 1 "AWK"
Set-TextCell 163 2 "05/18/2023"
Set-PlainCell 163 3 "Michael Marberry"
$ws.Cells.Item(163,4).Value2 = 1400
Set-TextCell 163 5 "`$142.35"
$ws.Cells.Item(163,6).Value2 = 199290
$ws.Cells.Item(163,7).Value2 = 61.59

# Row 164
Set-PlainCell 164 1 "AWK"
Set-TextCell 164 2 "05/02/2023"
Set-PlainCell 164 3 "Michael Marberry"
$ws.Cells.Item(164,4).Value2 = 675
Set-TextCell 164 5 "`$145.89"
$ws.Cells.Item(164,6).Value2 = 98475.74999999999
$ws.Cells.Item(164,7).Value2 = 119.89

# Row 165
Set-PlainCell 165 1 "AWK"
Set-TextCell 165 2 "03/01/2021"
Set-PlainCell 165 3 "MARTHA CLARK GOSS"
$ws.Cells.Item(165,4).Value2 = 61
Set-TextCell 165 5 "`$144.83"
$ws.Cells.Item(165,6).Value2 = 8834.630000000001
$ws.Cells.Item(165,7).Value2 = 0.19

# Row 166
Set-PlainCell 166 1 "AWK"
Set-TextCell 166 2 "05/13/2020"
Set-PlainCell 166 3 "Lloyd M Yates"
$ws.Cells.Item(166,4).Value2 = 2000
Set-TextCell 166 5 "`$115.95"
$ws.Cells.Item(166,6).Value2 = 231900
$ws.Cells.Item(166,7).Value2 = 48.78

# Row 167
Set-PlainCell 167 1 "AMGN"
Set-TextCell 167 2 "09/19/2023"
Set-PlainCell 167 3 "INC AMGEN"
$ws.Cells.Item(167,4).Value2 = 1764705
Set-TextCell 167 5 "`$17.00"
$ws.Cells.Item(167,6).Value2 = 29999985
$ws.Cells.Item(167,7).Value2 = 5.25
